$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 43 new blank rows before row 2 (shift existing data rows 2-121 down to 45-164)
$ws.Rows.Item(2).Resize(43).Insert()

# Copy formatting from the row immediately below (row 45, which retains the original row-2 style)
# so the newly inserted rows match the existing date/number formatting instead of Excel's
# "inherit from row above" default.
$ws.Range("A45").Copy()
$ws.Range("A2:A44").PasteSpecial(-4122)
$ws.Range("B45").Copy()
$ws.Range("B2:B44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New backward-extension data (43 quarters of real-time GDP revision history)
$newDates = @(30864,30956,31048,31138,31229,31321,31413,31503,31594,31686,31778,31868,31959,32051,32143,32234,32325,32417,32509,32599,32690,32782,32874,32964,33055,33147,33239,33329,33420,33512,33604,33695,33786,33878,33970,34060,34151,34243,34335,34425,34516,34608,34700)
$newValues = @(58.1640640589607,58.65126972831932,58.34028738617552,58.91042168010581,59.67751145739383,60.04032418989493,59.69824361353676,60.39277084432455,60.83851220139732,61.43974472954199,59.89519909689449,61.24278924618426,61.67816452518557,62.65257586390278,62.09280764804395,63.22271015783307,63.91723738862088,64.68432716590888,65.34775616248231,65.61727419234028,66.20814064241347,66.98559649777297,68.37465095934856,68.70636545763529,70.15761638763965,71.4015457562148,73.37110058979218,75.28639195525506,74.54262576899043,74.8124746887208,75.77352864378278,76.15256821218293,76.11904662965016,76.12584519098394,74.04422648832355,74.92863296232139,75.58528044338553,75.35079842925288,76.08801065791545,77.49682360591113,78.01155022377557,78.12627684164005,78.75139914350594)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newValues[$i]
}

# The pre-existing rows (now at 45-164) were re-dated: every date moves back by 45 days
# while the revision values themselves stay exactly as they were.
for ($r = 45; $r -le 164; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldDate = $cell.Value2()
    $cell.Value = $oldDate - 45
}
